$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Infrastructure")

# Set the "Alias" (row 14) and "Keystore password / key password" (row 15)
# values in column C to "chovanhan", matching the Username value used
# elsewhere in the sheet (for APK signing credentials).
$ws.Range("C14").Value = "chovanhan"
$ws.Range("C15").Value = "chovanhan"

# Update the active selection to C16, as recorded in the saved workbook view.
$ws.Range("C16").Select()
